$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2025
$ws.Range("I2").Value = 604.75
$ws.Range("J2").Value = 3057.9092
$ws.Range("K2").Value = 604.75
$ws.Range("L2").Value = 3057.9092
$ws.Range("M2").Value = -491.75
$ws.Range("N2").Value = -3283.9092
$ws.Range("H62").Value = 7919.1665
$ws.Range("J62").Value = 7927.4
$ws.Range("L62").Value = 7927.4
$ws.Range("N62").Value = -9175.4
$ws.Range("H64").Value = 5500
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 5500
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 5500
$ws.Range("N64").Value = -5996
$ws.Range("H65").Value = 7919.1665
$ws.Range("J65").Value = 7927.4
$ws.Range("L65").Value = 39637
$ws.Range("N65").Value = -45877
$ws.Range("H67").Value = 5500
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 5500
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 5500
$ws.Range("N67").Value = -7216
$ws.Range("H74").Value = 126146.75
$ws.Range("I74").Value = 186043.12
$ws.Range("J74").Value = 6354
$ws.Range("K74").Value = 186043.12
$ws.Range("L74").Value = 6354
$ws.Range("M74").Value = -185107.12
$ws.Range("N74").Value = -8226
$ws.Range("H77").Value = 126146.75
$ws.Range("I77").Value = 186043.12
$ws.Range("J77").Value = 6354
$ws.Range("K77").Value = 930215.6
$ws.Range("L77").Value = 31770
$ws.Range("M77").Value = -925535.6
$ws.Range("N77").Value = -41130
$ws.Range("H92").Value = 293.36365
$ws.Range("I92").Value = 338.6
$ws.Range("J92").Value = 255.66667
$ws.Range("K92").Value = 338.6
$ws.Range("L92").Value = 255.66667
$ws.Range("M92").Value = 909.4
$ws.Range("N92").Value = -2751.66667
$ws.Range("H113").Value = 2929.1667
$ws.Range("I113").Value = 2793.75
$ws.Range("K113").Value = 2793.75
$ws.Range("M113").Value = 460.25
$ws.Range("H141").Value = 3070.4285
$ws.Range("I141").Value = 3070.4285
$ws.Range("K141").Value = 9211.2855
$ws.Range("M141").Value = -4031.2855
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4952.6
$ws.Range("I36").Value = 3690.75
$ws.Range("K36").Value = 3690.75
$ws.Range("M36").Value = -3344.75
$ws.Range("H97").Value = 999.6667
$ws.Range("J97").Value = 999
$ws.Range("L97").Value = 999
$ws.Range("H132").Value = 2708.9092
$ws.Range("I132").Value = 2400
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 7200
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -4670
$ws.Range("N132").Value = -14808.5
$ws.Range("N97").Value = -1991

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 613.6667
$ws.Range("I22").Value = 485.375
$ws.Range("K22").Value = 485.375
$ws.Range("M22").Value = -312.375
$ws.Range("H86").Value = 694
$ws.Range("I86").Value = 694
$ws.Range("K86").Value = 694
$ws.Range("M86").Value = 429
$ws.Range("H89").Value = 694
$ws.Range("I89").Value = 694
$ws.Range("K89").Value = 3470
$ws.Range("M89").Value = 2146
$ws.Range("H94").Value = 1810.5714
$ws.Range("I94").Value = 1782.3334
$ws.Range("J94").Value = 1980
$ws.Range("K94").Value = 1782.3334
$ws.Range("L94").Value = 1980
$ws.Range("M94").Value = -1331.3334
$ws.Range("H99").Value = 2015.5385
$ws.Range("I99").Value = 2026
$ws.Range("K99").Value = 2026
$ws.Range("M99").Value = -528
$ws.Range("N94").Value = -2882

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1130.1818
$ws.Range("I16").Value = 899.875
$ws.Range("J16").Value = 1744.3334
$ws.Range("K16").Value = 899.875
$ws.Range("L16").Value = 1744.3334
$ws.Range("M16").Value = -612.875
$ws.Range("N16").Value = -2318.3334
$ws.Range("H58").Value = 2424.3572
$ws.Range("I58").Value = 2411.75
$ws.Range("K58").Value = 2411.75
$ws.Range("M58").Value = -2208.75
$ws.Range("H105").Value = 3915.7144
$ws.Range("I105").Value = 3371.2856
$ws.Range("K105").Value = 3371.2856
$ws.Range("M105").Value = -1624.2856
$ws.Range("H113").Value = 1130.1818
$ws.Range("I113").Value = 899.875
$ws.Range("J113").Value = 1744.3334
$ws.Range("K113").Value = 899.875
$ws.Range("L113").Value = 1744.3334
$ws.Range("M113").Value = 1270.125
$ws.Range("N113").Value = -6084.3334
$ws.Range("H136").Value = 2424.3572
$ws.Range("I136").Value = 2411.75
$ws.Range("K136").Value = 7235.25
$ws.Range("M136").Value = -4685.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 36910.5
$ws.Range("I62").Value = 38056.8
$ws.Range("K62").Value = 38056.8
$ws.Range("M62").Value = -37370.8
$ws.Range("H65").Value = 36910.5
$ws.Range("I65").Value = 38056.8
$ws.Range("K65").Value = 114170.4
$ws.Range("M65").Value = -110738.4
$ws.Range("H70").Value = 8134.0713
$ws.Range("I70").Value = 7989.6665
$ws.Range("J70").Value = 8394
$ws.Range("K70").Value = 7989.6665
$ws.Range("L70").Value = 8394
$ws.Range("M70").Value = -7719.6665
$ws.Range("N70").Value = -8934
$ws.Range("H73").Value = 8134.0713
$ws.Range("I73").Value = 7989.6665
$ws.Range("J73").Value = 8394
$ws.Range("K73").Value = 7989.6665
$ws.Range("L73").Value = 8394
$ws.Range("M73").Value = -7053.6665
$ws.Range("N73").Value = -10266
$ws.Range("H97").Value = 259.6
$ws.Range("I97").Value = 287
$ws.Range("K97").Value = 287
$ws.Range("M97").Value = 209

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6899.423
$ws.Range("I7").Value = 2200
$ws.Range("K7").Value = 2200
$ws.Range("M7").Value = -2088
$ws.Range("H16").Value = 796
$ws.Range("I16").Value = 550.44446
$ws.Range("J16").Value = 1348.5
$ws.Range("K16").Value = 550.44446
$ws.Range("L16").Value = 1348.5
$ws.Range("M16").Value = -380.44446
$ws.Range("N16").Value = -1688.5
$ws.Range("H55").Value = 241.8
$ws.Range("I55").Value = 197.94444
$ws.Range("K55").Value = 197.94444
$ws.Range("M55").Value = -24.94443999999999
$ws.Range("H68").Value = 2622.2222
$ws.Range("I68").Value = 1931.3334
$ws.Range("J68").Value = 2967.6667
$ws.Range("K68").Value = 1931.3334
$ws.Range("L68").Value = 2967.6667
$ws.Range("M68").Value = -1182.3334
$ws.Range("N68").Value = -4465.6667
$ws.Range("H71").Value = 2622.2222
$ws.Range("I71").Value = 1931.3334
$ws.Range("J71").Value = 2967.6667
$ws.Range("K71").Value = 9656.666999999999
$ws.Range("L71").Value = 14838.3335
$ws.Range("M71").Value = -5912.666999999999
$ws.Range("N71").Value = -22326.3335
$ws.Range("H93").Value = 1079.8
$ws.Range("I93").Value = 899.5
$ws.Range("K93").Value = 899.5
$ws.Range("H126").Value = 6899.423
$ws.Range("I126").Value = 2200
$ws.Range("K126").Value = 6600
$ws.Range("M126").Value = -4130
$ws.Range("H132").Value = 3104.1333
$ws.Range("I132").Value = 2543.375
$ws.Range("K132").Value = 7630.125
$ws.Range("M132").Value = -5100.125
$ws.Range("M93").Value = 348.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2242.2144
$ws.Range("I136").Value = 2030.0769
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6090.2307
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -3540.2307
$ws.Range("N136").Value = -20100
